$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date in A1 (45406 -> 45436)
$ws.Range("A1").Value = 45436

# Update price list values in column D
$ws.Range("D28").Value = 935
$ws.Range("D29").Value = 1020
$ws.Range("D30").Value = 1150
$ws.Range("D31").Value = 1350
$ws.Range("D32").Value = 1640
$ws.Range("D33").Value = 1790
$ws.Range("D34").Value = 1995
$ws.Range("D35").Value = 2210
$ws.Range("D36").Value = 2440
$ws.Range("D37").Value = 2996
$ws.Range("D38").Value = 3360
$ws.Range("D39").Value = 3750
$ws.Range("D40").Value = 4190
$ws.Range("D41").Value = 4550
$ws.Range("D42").Value = 4960
$ws.Range("D43").Value = 5270
$ws.Range("D49").Value = 2820
$ws.Range("D50").Value = 2992
$ws.Range("D51").Value = 3601
$ws.Range("D52").Value = 3850
$ws.Range("D53").Value = 3950
$ws.Range("D54").Value = 4260
$ws.Range("D55").Value = 5050
$ws.Range("D56").Value = 5670
$ws.Range("D57").Value = 6550
$ws.Range("D58").Value = 7250
$ws.Range("D59").Value = 7990
$ws.Range("D60").Value = 8920
$ws.Range("D61").Value = 9680
$ws.Range("D62").Value = 10370
$ws.Range("D68").Value = 1233
$ws.Range("D69").Value = 1284
$ws.Range("D70").Value = 1449
$ws.Range("D71").Value = 1650
$ws.Range("D72").Value = 1900
$ws.Range("D73").Value = 2130
$ws.Range("D74").Value = 2510
$ws.Range("D75").Value = 2670
$ws.Range("D76").Value = 3110
$ws.Range("D77").Value = 3630
$ws.Range("D78").Value = 4102
$ws.Range("D79").Value = 4650
$ws.Range("D80").Value = 5115
$ws.Range("D81").Value = 5310
$ws.Range("D82").Value = 5930
$ws.Range("D83").Value = 6496
$ws.Range("D89").Value = 556
$ws.Range("D90").Value = 591
$ws.Range("D91").Value = 635
$ws.Range("D92").Value = 765
$ws.Range("D93").Value = 951
$ws.Range("D94").Value = 1050
$ws.Range("D95").Value = 1080
$ws.Range("D96").Value = 1260
$ws.Range("D97").Value = 1470
$ws.Range("D98").Value = 1760
$ws.Range("D99").Value = 1930
$ws.Range("D100").Value = 2220
$ws.Range("D106").Value = 4790
$ws.Range("D107").Value = 5380
$ws.Range("D108").Value = 5670
$ws.Range("D109").Value = 6020
$ws.Range("D110").Value = 6350
$ws.Range("D111").Value = 7070
$ws.Range("D112").Value = 8080
$ws.Range("D113").Value = 9420
$ws.Range("D114").Value = 10400
$ws.Range("D115").Value = 11440
$ws.Range("D116").Value = 12580
$ws.Range("D117").Value = 13650
$ws.Range("D118").Value = 14550
$ws.Range("D119").Value = 15820
